$d = $word.ActiveDocument

# --- Fix 1: merge the "[Step 4]" run with the preceding space run -----------
# In the source document, the sentence "A pop-up appears ... record." is
# followed by two separate runs: a single space run and a "[Step 4]" run
# (same character formatting, different w:rsidR). Editing only the
# "[Step 4]" run (not touching the space run) and replacing its text with
# itself causes the two adjacent, identically-formatted runs to coalesce
# into one run, exactly mirroring what Word does after a Find/Replace or
# retype over that span.
$rStep4 = $d.Content
$rStep4.Find.Execute("[Step 4]")
$narrow = $d.Range($rStep4.Start, $rStep4.End)
$narrow.Find.Execute("[Step 4]", $false, $false, $false, $false, $false, $true, 1, $false, "[Step 4]", 2) | Out-Null

# --- Fix 2: correct the "f0orm" typo to "form" -------------------------------
# Locate the stray "0" inside "f0orm" and remove just that character.
$rTypo = $d.Content
$rTypo.Find.Execute("f0orm")
$zeroStart = $rTypo.Start + 1
$zeroRange = $d.Range($zeroStart, $zeroStart + 1)
$zeroRange.Find.Execute("0", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# The edit point (collapsed range right after "f") is where Word's automatic
# "_GoBack" bookmark (tracking the most recent edit) now belongs. Move the
# existing bookmark there.
$editPoint = $zeroStart
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()
$newGoBackRange = $d.Range($editPoint, $editPoint)
$d.Bookmarks.Add("_GoBack", $newGoBackRange)
